# Add a "Group" column to the site report template.
# This mirrors a user inserting a new column before column E (which held
# "Count Code"/"Count"), typing the header "Group" in the new E3 cell, and
# leaving the new column selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E. Excel automatically shifts every
# column from E onward to the right by one and copies the formatting of
# the column to the left (D) into the newly inserted column.
$ws.Columns("E:E").Insert()

# Match the width of the new column to its neighbour (column D), so the
# <col> run for the style-7 formatted block effectively spans columns D:E.
$ws.Columns("E").ColumnWidth = $ws.Columns("D").ColumnWidth

# Give the new column its header text.
$ws.Range("E3").Value = "Group"

# Leave the new header cell selected, matching the saved cursor position.
$ws.Range("E3").Select() | Out-Null
